$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The currency ("Moeda") formatting previously applied to the Saldo column
# is no longer needed - normalize those cells (and the header cell above
# them) back to the plain "Normal" style, re-adding just the bottom border
# under the header row to match the rest of row 1.
$ws.Range("C1:C20").Style = "Normal"
$ws.Range("A1:C1").Borders.Item(9).LineStyle = 1

# Append the applied-filters note as a new row below the data.
$ws.Range("A22").Value = "Filtros aplicados:
DataFim é (Em branco)
nr_saldo_disponivel não é 0
Posição é Posição D-1
CARTEIRA não está em branco
DataFim é (Em branco)
NR_CONTA não está em branco
TIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"

# Reflect the new active cell position left behind after the edit.
$ws.Range("G3").Select()
